# Update imputed values for specific cells in columns A and C
# (Random Forest algorithm result data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4" = -21.48709999999999
    "A6" = -22.54160000000002
    "A7" = -19.36759999999999
    "C7" = -12.3717
    "A8" = -22.39530000000002
    "C11" = -11.7229
    "C12" = -10.5805
    "C15" = -14.75869999999999
    "A16" = -21.53659999999998
    "A20" = -19.13399999999999
    "C20" = -12.0704
    "A21" = -19.94059999999998
    "C21" = -11.849
    "C22" = -12.52770000000001
    "C23" = -12.29710000000001
    "A28" = -21.7949
    "A29" = -21.53329999999998
    "C29" = -11.44650000000001
    "A30" = -21.5461
    "A32" = -21.19979999999999
    "C34" = -11.16280000000001
    "A40" = -20.20889999999999
    "C42" = -11.9234
    "C43" = -13.47959999999998
    "C44" = -13.9323
    "C45" = -14.09609999999999
    "A46" = -21.9216
    "C46" = -12.86469999999999
    "C50" = -14.16459999999999
    "A51" = -21.61559999999999
    "C51" = -10.9406
    "A52" = -22.30640000000001
    "A57" = -22.0308
    "C57" = -12.58549999999999
    "A59" = -22.30390000000002
    "A62" = -22.30540000000001
    "C65" = -13.26009999999999
    "A66" = -21.61139999999999
    "C66" = -11.04570000000001
    "C67" = -11.1234
    "A73" = -20.68099999999999
    "A74" = -21.82919999999999
    "A77" = -19.98489999999999
    "C79" = -11.7492
    "C84" = -12.55609999999999
    "C87" = -13.89879999999999
    "A92" = -21.84329999999998
    "C92" = -11.35980000000001
    "C97" = -11.12500000000001
    "A100" = -22.2711
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
